$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# --- Simple price/value-only updates (column D, and a couple of E cells) ---
Set-TextValue "D2" "247.80"
Set-TextValue "D4" "5.457"
Set-TextValue "D5" "0.05664"
Set-TextValue "D7" "0.8005"
Set-TextValue "D8" "1.038"

Set-TextValue "D18" "0.006363"
Set-TextValue "D19" "0.005017"
Set-TextValue "D20" "0.001047"
Set-TextValue "D22" "0.0003201"
Set-TextValue "D23" "3.824"
Set-TextValue "D24" "6.427"
Set-TextValue "D25" "2.089"
Set-TextValue "D26" "0.3279"
Set-TextValue "D27" "0.1317"

Set-TextValue "D40" "0.04086"
Set-TextValue "D41" "0.006951"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

Set-TextValue "D44" "0.008950"
Set-TextValue "D45" "0.00005848"
Set-TextValue "D47" "0.7854"
Set-TextValue "D48" "0.01174"

# --- Rows 9-17: the coin ranking list shifted up by one position (row 9's
#     "One" entry dropped off the top and reappears at the bottom, rank 17) ---
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1443"
$ws.Range("E9").Value = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07207"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.03148"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.02939"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09286"
$ws.Range("E13").Value = "12BitMartTokenBMX"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001661"
$ws.Range("E14").Value = "13BitForexTokenBF"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.211"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04727"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D17" "0.0005891"
$ws.Range("E17").Value = "16OneONE"

# --- Rows 42-43: BKEXToken and CEJI swapped order ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003501"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1040"
$ws.Range("E43").Value = "42BKEXTokenBKK"
